$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the source workbook, where these columns are stored as inline strings)
# so that e.g. "0.630" is not silently normalized to the number 0.63.
$textForceRefs = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D26",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D47",
    "D49",
    "D51",
)
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data
$ws.Range("D2").Value = '42.643.91'
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").Value = '2.252.38'
$ws.Range("E3").Value = '  +0.48%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '246.27'
$ws.Range("E5").Value = '  -0.28%  '

$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").Value = '76.63'
$ws.Range("E7").Value = '  +2.15%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '0.630'
$ws.Range("E9").Value = '  +0.13%  '

$ws.Range("D10").Value = '45.11'
$ws.Range("E10").Value = '  +12.20%  '

$ws.Range("D11").Value = '0.0955'
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").Value = '7.35'
$ws.Range("E12").Value = '  +2.37%  '

$ws.Range("E13").Value = '  -0.19%  '

$ws.Range("D14").Value = '14.80'
$ws.Range("E14").Value = '  -0.41%  '

$ws.Range("D15").Value = '0.867'
$ws.Range("E15").Value = '  +0.47%  '

$ws.Range("D16").Value = '2.237.27'
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").Value = '42.430.44'

$ws.Range("E18").Value = '  +3.89%  '

$ws.Range("D19").Value = '6.22'
$ws.Range("E19").Value = '  +1.30%  '

$ws.Range("D20").Value = '72.19'
$ws.Range("E20").Value = '  +0.77%  '

$ws.Range("D21").Value = '11.09'
$ws.Range("E21").Value = '  +55.22%  '

$ws.Range("E22").Value = '  +1.29%  '

$ws.Range("D23").Value = '232.57'
$ws.Range("E23").Value = '  +0.66%  '

$ws.Range("D24").Value = '11.94'
$ws.Range("E24").Value = '  +4.22%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").Value = '3.63'
$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("E28").Value = '  +5.44%  '

$ws.Range("D29").Value = '167.45'
$ws.Range("E29").Value = '  -0.91%  '

$ws.Range("D30").Value = '20.75'
$ws.Range("E30").Value = '  +1.07%  '

$ws.Range("D31").Value = '0.0826'
$ws.Range("E31").Value = '  -2.14%  '

$ws.Range("D32").Value = '32.50'
$ws.Range("E32").Value = '  -3.72%  '

$ws.Range("D33").Value = '5.71'
$ws.Range("E33").Value = '  +17.00%  '

$ws.Range("D34").Value = '0.121'
$ws.Range("E34").Value = '  -0.50%  '

$ws.Range("D35").Value = '0.126'
$ws.Range("E35").Value = '  -0.57%  '

$ws.Range("D36").Value = '4.74'
$ws.Range("E36").Value = '  +4.11%  '

$ws.Range("D37").Value = '0.0318'
$ws.Range("E37").Value = '  +5.99%  '

$ws.Range("D38").Value = '14.38'
$ws.Range("E38").Value = '  +6.88%  '

$ws.Range("E39").Value = '  +0.52%  '

$ws.Range("D40").Value = '5.81'
$ws.Range("E40").Value = '  -1.89%  '

$ws.Range("D41").Value = '64.32'
$ws.Range("E41").Value = '  +6.32%  '

$ws.Range("D42").Value = '0.204'
$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("D43").Value = '108.51'
$ws.Range("E43").Value = '  -3.37%  '

$ws.Range("D44").Value = '8.95'
$ws.Range("E44").Value = '  +1.75%  '

$ws.Range("E45").Value = '  +2.21%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").Value = '2.41'
$ws.Range("E47").Value = '  +7.01%  '

$ws.Range("E48").Value = '  +1.05%  '

$ws.Range("D49").Value = '1.19'
$ws.Range("E49").Value = '  +1.64%  '

$ws.Range("E50").Value = '  -1.26%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = '0.426'
$ws.Range("E51").Value = '  +11.06%  '

# Restore default (General) style on cells that were temporarily switched to
# text format, now that the text type itself has "stuck" to the cell.
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
